$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new data rows (obstacles / coin events) below the existing header row
$ws.Range("A2").Value = "Fri Dec  8 08_33_51 2023"
$ws.Range("B2").Value = "loclexyz99"
$ws.Range("C2").Value = -1

$ws.Range("A3").Value = "Fri Dec  8 09_18_54 2023"
$ws.Range("B3").Value = "loclexyz99"
$ws.Range("C3").Value = -3

# Move the active selection, matching the saved view state
$ws.Range("F6").Select() | Out-Null
